$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H13").Value = 10000.5
$ws.Range("I13").Value = 0
$ws.Range("K13").Value = 0
$ws.Range("M13").Value = ""
$ws.Range("H28").Value = 1378.9048
$ws.Range("I28").Value = 1155.3158
$ws.Range("K28").Value = 1155.3158
$ws.Range("M28").Value = -670.3158000000001
$ws.Range("H32").Value = 10000
$ws.Range("I32").Value = 10000
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 10000
$ws.Range("L32").Value = 0
$ws.Range("M32").Value = -9674
$ws.Range("N32").Value = ""
$ws.Range("H70").Value = 2319.2307
$ws.Range("I70").Value = 1672.2222
$ws.Range("J70").Value = 3775
$ws.Range("K70").Value = 5016.6666
$ws.Range("L70").Value = 11325
$ws.Range("M70").Value = -4746.6666
$ws.Range("N70").Value = -11865
$ws.Range("H73").Value = 2319.2307
$ws.Range("I73").Value = 1672.2222
$ws.Range("J73").Value = 3775
$ws.Range("K73").Value = 5016.6666
$ws.Range("L73").Value = 11325
$ws.Range("M73").Value = -4080.6666
$ws.Range("N73").Value = -13197
$ws.Range("H82").Value = 0
$ws.Range("I82").Value = 0
$ws.Range("K82").Value = 0
$ws.Range("M82").Value = ""
$ws.Range("H85").Value = 0
$ws.Range("I85").Value = 0
$ws.Range("K85").Value = 0
$ws.Range("M85").Value = ""
$ws.Range("H137").Value = 33334922
$ws.Range("I137").Value = 66667844
$ws.Range("K137").Value = 200003532
$ws.Range("M137").Value = -200000982
$ws.Range("H141").Value = 1000
$ws.Range("I141").Value = 0
$ws.Range("K141").Value = 0
$ws.Range("M141").Value = ""

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H6").Value = 60000000
$ws.Range("I6").Value = 60000000
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 60000000
$ws.Range("L6").Value = 0
$ws.Range("M6").Value = -59999827
$ws.Range("N6").Value = ""
$ws.Range("H14").Value = 868.6667
$ws.Range("I14").Value = 868.6667
$ws.Range("J14").Value = 0
$ws.Range("K14").Value = 868.6667
$ws.Range("L14").Value = 0
$ws.Range("M14").Value = -693.6667
$ws.Range("N14").Value = ""
$ws.Range("H61").Value = 2000
$ws.Range("I61").Value = 2000
$ws.Range("K61").Value = 2000
$ws.Range("M61").Value = -1788
$ws.Range("H76").Value = 48530.668
$ws.Range("I76").Value = 30000
$ws.Range("J76").Value = 52236.8
$ws.Range("K76").Value = 30000
$ws.Range("L76").Value = 52236.8
$ws.Range("M76").Value = -29662
$ws.Range("N76").Value = -52912.8
$ws.Range("H79").Value = 48530.668
$ws.Range("I79").Value = 30000
$ws.Range("J79").Value = 52236.8
$ws.Range("K79").Value = 30000
$ws.Range("L79").Value = 52236.8
$ws.Range("M79").Value = -28830
$ws.Range("N79").Value = -54576.8
$ws.Range("H110").Value = 2792.2
$ws.Range("I110").Value = 3755.5
$ws.Range("J110").Value = 2150
$ws.Range("K110").Value = 3755.5
$ws.Range("L110").Value = 2150
$ws.Range("M110").Value = -1710.5
$ws.Range("N110").Value = -6240
$ws.Range("H136").Value = 2000
$ws.Range("I136").Value = 2000
$ws.Range("K136").Value = 6000
$ws.Range("M136").Value = -3450

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 17891.766
$ws.Range("I86").Value = 20043.154
$ws.Range("K86").Value = 20043.154
$ws.Range("M86").Value = -18920.154
$ws.Range("H89").Value = 17891.766
$ws.Range("I89").Value = 20043.154
$ws.Range("K89").Value = 100215.77
$ws.Range("M89").Value = -94599.76999999999
$ws.Range("H99").Value = 3250
$ws.Range("J99").Value = 5000
$ws.Range("L99").Value = 5000
$ws.Range("N99").Value = -7996
$ws.Range("H105").Value = 4062.5715
$ws.Range("I105").Value = 4062.5715
$ws.Range("K105").Value = 4062.5715
$ws.Range("M105").Value = -2315.5715
$ws.Range("H134").Value = 1828.5625
$ws.Range("I134").Value = 1825.5
$ws.Range("J134").Value = 1850
$ws.Range("K134").Value = 5476.5
$ws.Range("L134").Value = 5550
$ws.Range("M134").Value = -2941.5
$ws.Range("N134").Value = -10620

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H11").Value = 700
$ws.Range("J11").Value = 0
$ws.Range("L11").Value = 0
$ws.Range("N11").Value = ""

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 6000157.5
$ws.Range("I4").Value = 197
$ws.Range("K4").Value = 591
$ws.Range("M4").Value = -479
$ws.Range("H34").Value = 81612.62
$ws.Range("I34").Value = 40
$ws.Range("J34").Value = 96444
$ws.Range("K34").Value = 120
$ws.Range("L34").Value = 289332
$ws.Range("M34").Value = -36
$ws.Range("N34").Value = -289500
$ws.Range("H39").Value = 6637.4375
$ws.Range("J39").Value = 6579.933
$ws.Range("L39").Value = 19739.799
$ws.Range("N39").Value = -20327.799
$ws.Range("H40").Value = 106.833336
$ws.Range("I40").Value = 48.2
$ws.Range("K40").Value = 192.8
$ws.Range("M40").Value = -123.8
$ws.Range("H62").Value = 10399.4
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 10399.4
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 31198.2
$ws.Range("M62").Value = ""
$ws.Range("N62").Value = -32570.2
$ws.Range("H65").Value = 10399.4
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 10399.4
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 93594.59999999999
$ws.Range("M65").Value = ""
$ws.Range("N65").Value = -100458.6
$ws.Range("H86").Value = 566.3333
$ws.Range("J86").Value = 223.25
$ws.Range("L86").Value = 669.75
$ws.Range("N86").Value = -3041.75
$ws.Range("H89").Value = 566.3333
$ws.Range("J89").Value = 223.25
$ws.Range("L89").Value = 2009.25
$ws.Range("N89").Value = -13865.25
$ws.Range("H113").Value = 2358.2856
$ws.Range("I113").Value = 2358.2856
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 7074.8568
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = -4904.8568
$ws.Range("N113").Value = ""
$ws.Range("H133").Value = 7366
$ws.Range("I133").Value = 7366
$ws.Range("K133").Value = 22098
$ws.Range("M133").Value = -17038

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 20836256
$ws.Range("I132").Value = 2540
$ws.Range("K132").Value = 7620
$ws.Range("M132").Value = -5090

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 3668.9412
$ws.Range("I46").Value = 1574.2858
$ws.Range("K46").Value = 1574.2858
$ws.Range("M46").Value = -1386.2858
$ws.Range("H93").Value = 720.875
$ws.Range("I93").Value = 668.1429000000001
$ws.Range("K93").Value = 668.1429000000001
$ws.Range("M93").Value = 579.8570999999999
$ws.Range("H95").Value = 40000
$ws.Range("J95").Value = 40000
$ws.Range("L95").Value = 40000
$ws.Range("N95").Value = -45492
$ws.Range("H122").Value = 1796.2
$ws.Range("I122").Value = 1796.2
$ws.Range("K122").Value = 5388.6
$ws.Range("M122").Value = -2938.6
$ws.Range("H132").Value = 4149.75
$ws.Range("J132").Value = 3533.3333
$ws.Range("L132").Value = 10599.9999
$ws.Range("N132").Value = -15659.9999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H11").Value = 4233.385
$ws.Range("I11").Value = 4233.385
$ws.Range("K11").Value = 4233.385
$ws.Range("M11").Value = -4091.385
$ws.Range("H14").Value = 1548.75
$ws.Range("I14").Value = 2247.5
$ws.Range("K14").Value = 2247.5
$ws.Range("M14").Value = -2079.5
$ws.Range("H32").Value = 3341504.2
$ws.Range("I32").Value = 6668675.5
$ws.Range("J32").Value = 14333.333
$ws.Range("K32").Value = 6668675.5
$ws.Range("L32").Value = 14333.333
$ws.Range("M32").Value = -6668358.5
$ws.Range("N32").Value = -14967.333
$ws.Range("H98").Value = 37999.5
$ws.Range("J98").Value = 35999
$ws.Range("L98").Value = 35999
$ws.Range("N98").Value = -41989
$ws.Range("H127").Value = 0
$ws.Range("J127").Value = 0
$ws.Range("L127").Value = 0
$ws.Range("N127").Value = ""
